# Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Ajo
# A new weekly observation row is inserted at row 508 (shifting the
# existing rows 508:561 down to 509:562), extending the data range
# from A1:R561 to A1:R562.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 508, pushing everything
# from 508 downward by one row.
$ws.Rows.Item(508).Insert()

# Re-apply the date number format used throughout column D to the
# newly inserted cell (copy it from the row right below, which still
# holds the original formatting).
$ws.Cells.Item(508, 4).NumberFormat = $ws.Cells.Item(509, 4).NumberFormat

# Populate the new row with the new weekly record.
$ws.Cells.Item(508, 1).Value  = 6
$ws.Cells.Item(508, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(508, 3).Value  = "Metropolitana"
$ws.Cells.Item(508, 4).Value  = (Get-Date -Year 2022 -Month 1 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(508, 5).Value  = 13
$ws.Cells.Item(508, 6).Value  = 100112003
$ws.Cells.Item(508, 7).Value  = "Ajo"
$ws.Cells.Item(508, 8).Value  = "Chino"
$ws.Cells.Item(508, 9).Value  = "Primera"
$ws.Cells.Item(508, 10).Value = 1400
$ws.Cells.Item(508, 11).Value = 16500
$ws.Cells.Item(508, 12).Value = 17000
$ws.Cells.Item(508, 13).Value = 16679
$ws.Cells.Item(508, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(508, 15).Value = "China"
$ws.Cells.Item(508, 16).Value = 1668
$ws.Cells.Item(508, 17).Value = 10
$ws.Cells.Item(508, 18).Value = "Hortaliza"
